# Progress update as of 04-Nov-2025.
# For every data row (3..35) on the "Training Dashboard" sheet:
#   - column H ("PERIOD TO EXPIRE") ticks down by 1 day
#   - column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025
#
# The new "LAST UPDATE" text is forced to text with a leading apostrophe so
# Excel stores it as the literal string "04-Nov-2025" instead of silently
# recognising the pattern and auto-converting it to a date serial number
# (which is what a bare `Value = "04-Nov-2025"` assignment would do).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 35; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)   # column H
    $updateCell = $ws.Cells.Item($row, 9)   # column I

    $current = $periodCell.Value2
    $periodCell.Value = $current - 1

    $updateCell.Value = "'04-Nov-2025"
}
